$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the former last row (342), shifting it down to row 346
$ws.Rows.Item(342).Resize(4).Insert()

# Update existing rows 318-341 with revised daily price data
$ws.Cells.Item(318, 4).Value = 44578
$ws.Cells.Item(318, 9).Value = 'Extra'
$ws.Cells.Item(318, 10).Value = 160
$ws.Cells.Item(318, 11).Value = 2800
$ws.Cells.Item(318, 12).Value = 3000
$ws.Cells.Item(318, 13).Value = 2900
$ws.Cells.Item(318, 14).Value = '$/unidad'
$ws.Cells.Item(318, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(318, 16).Value = 2900

$ws.Cells.Item(319, 4).Value = 44578
$ws.Cells.Item(319, 9).Value = 'Primera'
$ws.Cells.Item(319, 10).Value = 250
$ws.Cells.Item(319, 11).Value = 2300
$ws.Cells.Item(319, 12).Value = 2500
$ws.Cells.Item(319, 13).Value = 2400
$ws.Cells.Item(319, 14).Value = '$/unidad'
$ws.Cells.Item(319, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(319, 16).Value = 2400

$ws.Cells.Item(320, 4).Value = 44578
$ws.Cells.Item(320, 9).Value = 'Segunda'
$ws.Cells.Item(320, 10).Value = 97
$ws.Cells.Item(320, 11).Value = 1800
$ws.Cells.Item(320, 12).Value = 2100
$ws.Cells.Item(320, 13).Value = 1948
$ws.Cells.Item(320, 14).Value = '$/unidad'
$ws.Cells.Item(320, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(320, 16).Value = 1948

$ws.Cells.Item(321, 4).Value = 44578
$ws.Cells.Item(321, 9).Value = 'Tercera'
$ws.Cells.Item(321, 10).Value = 43
$ws.Cells.Item(321, 11).Value = 1200
$ws.Cells.Item(321, 12).Value = 1500
$ws.Cells.Item(321, 13).Value = 1353
$ws.Cells.Item(321, 14).Value = '$/unidad'
$ws.Cells.Item(321, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(321, 16).Value = 1353

$ws.Cells.Item(322, 4).Value = 44490
$ws.Cells.Item(322, 9).Value = 'Primera'
$ws.Cells.Item(322, 10).Value = 250
$ws.Cells.Item(322, 11).Value = 800
$ws.Cells.Item(322, 12).Value = 1000
$ws.Cells.Item(322, 13).Value = 900
$ws.Cells.Item(322, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(322, 15).Value = 'Perú'
$ws.Cells.Item(322, 16).Value = 900

$ws.Cells.Item(323, 4).Value = 44491
$ws.Cells.Item(323, 9).Value = 'Primera'
$ws.Cells.Item(323, 10).Value = 250
$ws.Cells.Item(323, 11).Value = 800
$ws.Cells.Item(323, 12).Value = 900
$ws.Cells.Item(323, 13).Value = 850
$ws.Cells.Item(323, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(323, 15).Value = 'Perú'
$ws.Cells.Item(323, 16).Value = 850

$ws.Cells.Item(324, 4).Value = 44266
$ws.Cells.Item(324, 9).Value = 'Extra'
$ws.Cells.Item(324, 10).Value = 250
$ws.Cells.Item(324, 11).Value = 2800
$ws.Cells.Item(324, 12).Value = 2800
$ws.Cells.Item(324, 13).Value = 2800
$ws.Cells.Item(324, 14).Value = '$/unidad'
$ws.Cells.Item(324, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(324, 16).Value = 2800

$ws.Cells.Item(325, 4).Value = 44266
$ws.Cells.Item(325, 9).Value = 'Extra'
$ws.Cells.Item(325, 10).Value = 250
$ws.Cells.Item(325, 11).Value = 2800
$ws.Cells.Item(325, 12).Value = 2800
$ws.Cells.Item(325, 13).Value = 2800
$ws.Cells.Item(325, 14).Value = '$/unidad'
$ws.Cells.Item(325, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(325, 16).Value = 2800

$ws.Cells.Item(326, 4).Value = 44266
$ws.Cells.Item(326, 9).Value = 'Primera'
$ws.Cells.Item(326, 10).Value = 340
$ws.Cells.Item(326, 11).Value = 2500
$ws.Cells.Item(326, 12).Value = 2500
$ws.Cells.Item(326, 13).Value = 2500
$ws.Cells.Item(326, 14).Value = '$/unidad'
$ws.Cells.Item(326, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(326, 16).Value = 2500

$ws.Cells.Item(327, 4).Value = 44266
$ws.Cells.Item(327, 9).Value = 'Primera'
$ws.Cells.Item(327, 10).Value = 340
$ws.Cells.Item(327, 11).Value = 2500
$ws.Cells.Item(327, 12).Value = 2500
$ws.Cells.Item(327, 13).Value = 2500
$ws.Cells.Item(327, 14).Value = '$/unidad'
$ws.Cells.Item(327, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(327, 16).Value = 2500

$ws.Cells.Item(328, 4).Value = 44266
$ws.Cells.Item(328, 9).Value = 'Segunda'
$ws.Cells.Item(328, 10).Value = 160
$ws.Cells.Item(328, 11).Value = 2300
$ws.Cells.Item(328, 12).Value = 2300
$ws.Cells.Item(328, 13).Value = 2300
$ws.Cells.Item(328, 14).Value = '$/unidad'
$ws.Cells.Item(328, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(328, 16).Value = 2300

$ws.Cells.Item(329, 4).Value = 44533
$ws.Cells.Item(329, 9).Value = 'Extra'
$ws.Cells.Item(329, 10).Value = 79
$ws.Cells.Item(329, 11).Value = 4300
$ws.Cells.Item(329, 12).Value = 4500
$ws.Cells.Item(329, 13).Value = 4399
$ws.Cells.Item(329, 14).Value = '$/unidad'
$ws.Cells.Item(329, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(329, 16).Value = 4399

$ws.Cells.Item(330, 4).Value = 44533
$ws.Cells.Item(330, 9).Value = 'Primera'
$ws.Cells.Item(330, 10).Value = 160
$ws.Cells.Item(330, 11).Value = 3600
$ws.Cells.Item(330, 12).Value = 4000
$ws.Cells.Item(330, 13).Value = 3800
$ws.Cells.Item(330, 14).Value = '$/unidad'
$ws.Cells.Item(330, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(330, 16).Value = 3800

$ws.Cells.Item(331, 4).Value = 44533
$ws.Cells.Item(331, 9).Value = 'Segunda'
$ws.Cells.Item(331, 10).Value = 52
$ws.Cells.Item(331, 11).Value = 3000
$ws.Cells.Item(331, 12).Value = 3400
$ws.Cells.Item(331, 13).Value = 3200
$ws.Cells.Item(331, 14).Value = '$/unidad'
$ws.Cells.Item(331, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(331, 16).Value = 3200

$ws.Cells.Item(332, 4).Value = 44264
$ws.Cells.Item(332, 9).Value = 'Extra'
$ws.Cells.Item(332, 10).Value = 160
$ws.Cells.Item(332, 11).Value = 2800
$ws.Cells.Item(332, 12).Value = 2800
$ws.Cells.Item(332, 13).Value = 2800
$ws.Cells.Item(332, 14).Value = '$/unidad'
$ws.Cells.Item(332, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(332, 16).Value = 2800

$ws.Cells.Item(333, 4).Value = 44264
$ws.Cells.Item(333, 9).Value = 'Extra'
$ws.Cells.Item(333, 10).Value = 250
$ws.Cells.Item(333, 11).Value = 2800
$ws.Cells.Item(333, 12).Value = 2800
$ws.Cells.Item(333, 13).Value = 2800
$ws.Cells.Item(333, 14).Value = '$/unidad'
$ws.Cells.Item(333, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(333, 16).Value = 2800

$ws.Cells.Item(334, 4).Value = 44264
$ws.Cells.Item(334, 9).Value = 'Primera'
$ws.Cells.Item(334, 10).Value = 250
$ws.Cells.Item(334, 11).Value = 2500
$ws.Cells.Item(334, 12).Value = 2500
$ws.Cells.Item(334, 13).Value = 2500
$ws.Cells.Item(334, 14).Value = '$/unidad'
$ws.Cells.Item(334, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(334, 16).Value = 2500

$ws.Cells.Item(335, 4).Value = 44264
$ws.Cells.Item(335, 9).Value = 'Primera'
$ws.Cells.Item(335, 10).Value = 340
$ws.Cells.Item(335, 11).Value = 2500
$ws.Cells.Item(335, 12).Value = 2500
$ws.Cells.Item(335, 13).Value = 2500
$ws.Cells.Item(335, 14).Value = '$/unidad'
$ws.Cells.Item(335, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(335, 16).Value = 2500

$ws.Cells.Item(336, 4).Value = 44264
$ws.Cells.Item(336, 9).Value = 'Segunda'
$ws.Cells.Item(336, 10).Value = 160
$ws.Cells.Item(336, 11).Value = 2200
$ws.Cells.Item(336, 12).Value = 2200
$ws.Cells.Item(336, 13).Value = 2200
$ws.Cells.Item(336, 14).Value = '$/unidad'
$ws.Cells.Item(336, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(336, 16).Value = 2200

$ws.Cells.Item(337, 4).Value = 44571
$ws.Cells.Item(337, 9).Value = 'Extra'
$ws.Cells.Item(337, 10).Value = 160
$ws.Cells.Item(337, 11).Value = 2800
$ws.Cells.Item(337, 12).Value = 3000
$ws.Cells.Item(337, 13).Value = 2900
$ws.Cells.Item(337, 14).Value = '$/unidad'
$ws.Cells.Item(337, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(337, 16).Value = 2900

$ws.Cells.Item(338, 4).Value = 44571
$ws.Cells.Item(338, 9).Value = 'Primera'
$ws.Cells.Item(338, 10).Value = 250
$ws.Cells.Item(338, 11).Value = 2400
$ws.Cells.Item(338, 12).Value = 2600
$ws.Cells.Item(338, 13).Value = 2500
$ws.Cells.Item(338, 14).Value = '$/unidad'
$ws.Cells.Item(338, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(338, 16).Value = 2500

$ws.Cells.Item(339, 4).Value = 44571
$ws.Cells.Item(339, 9).Value = 'Segunda'
$ws.Cells.Item(339, 10).Value = 106
$ws.Cells.Item(339, 11).Value = 2000
$ws.Cells.Item(339, 12).Value = 2200
$ws.Cells.Item(339, 13).Value = 2100
$ws.Cells.Item(339, 14).Value = '$/unidad'
$ws.Cells.Item(339, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(339, 16).Value = 2100

$ws.Cells.Item(340, 4).Value = 44571
$ws.Cells.Item(340, 9).Value = 'Tercera'
$ws.Cells.Item(340, 10).Value = 61
$ws.Cells.Item(340, 11).Value = 1600
$ws.Cells.Item(340, 12).Value = 1800
$ws.Cells.Item(340, 13).Value = 1702
$ws.Cells.Item(340, 14).Value = '$/unidad'
$ws.Cells.Item(340, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(340, 16).Value = 1702

$ws.Cells.Item(341, 4).Value = 44279
$ws.Cells.Item(341, 9).Value = 'Extra'
$ws.Cells.Item(341, 10).Value = 250
$ws.Cells.Item(341, 11).Value = 2500
$ws.Cells.Item(341, 12).Value = 2500
$ws.Cells.Item(341, 13).Value = 2500
$ws.Cells.Item(341, 14).Value = '$/unidad'
$ws.Cells.Item(341, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(341, 16).Value = 2500

# Fill in the newly inserted rows 342-345 with their full data
$ws.Cells.Item(342, 1).Value = 9
$ws.Cells.Item(342, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(342, 3).Value = 'Metropolitana'
$ws.Cells.Item(342, 5).Value = 13
$ws.Cells.Item(342, 6).Value = 100112028
$ws.Cells.Item(342, 7).Value = 'Sandia'
$ws.Cells.Item(342, 8).Value = 'Sin especificar'
$ws.Cells.Item(342, 17).Value = 1
$ws.Cells.Item(342, 18).Value = 'Hortaliza'
$ws.Cells.Item(342, 4).Value = 44279
$ws.Cells.Item(342, 9).Value = 'Primera'
$ws.Cells.Item(342, 10).Value = 340
$ws.Cells.Item(342, 11).Value = 2000
$ws.Cells.Item(342, 12).Value = 2000
$ws.Cells.Item(342, 13).Value = 2000
$ws.Cells.Item(342, 14).Value = '$/unidad'
$ws.Cells.Item(342, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(342, 16).Value = 2000

$ws.Cells.Item(343, 1).Value = 9
$ws.Cells.Item(343, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(343, 3).Value = 'Metropolitana'
$ws.Cells.Item(343, 5).Value = 13
$ws.Cells.Item(343, 6).Value = 100112028
$ws.Cells.Item(343, 7).Value = 'Sandia'
$ws.Cells.Item(343, 8).Value = 'Sin especificar'
$ws.Cells.Item(343, 17).Value = 1
$ws.Cells.Item(343, 18).Value = 'Hortaliza'
$ws.Cells.Item(343, 4).Value = 44279
$ws.Cells.Item(343, 9).Value = 'Segunda'
$ws.Cells.Item(343, 10).Value = 160
$ws.Cells.Item(343, 11).Value = 1600
$ws.Cells.Item(343, 12).Value = 1600
$ws.Cells.Item(343, 13).Value = 1600
$ws.Cells.Item(343, 14).Value = '$/unidad'
$ws.Cells.Item(343, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(343, 16).Value = 1600

$ws.Cells.Item(344, 1).Value = 9
$ws.Cells.Item(344, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(344, 3).Value = 'Metropolitana'
$ws.Cells.Item(344, 5).Value = 13
$ws.Cells.Item(344, 6).Value = 100112028
$ws.Cells.Item(344, 7).Value = 'Sandia'
$ws.Cells.Item(344, 8).Value = 'Sin especificar'
$ws.Cells.Item(344, 17).Value = 1
$ws.Cells.Item(344, 18).Value = 'Hortaliza'
$ws.Cells.Item(344, 4).Value = 44525
$ws.Cells.Item(344, 9).Value = 'Primera'
$ws.Cells.Item(344, 10).Value = 250
$ws.Cells.Item(344, 11).Value = 700
$ws.Cells.Item(344, 12).Value = 800
$ws.Cells.Item(344, 13).Value = 750
$ws.Cells.Item(344, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(344, 15).Value = 'Perú'
$ws.Cells.Item(344, 16).Value = 750

$ws.Cells.Item(345, 1).Value = 9
$ws.Cells.Item(345, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(345, 3).Value = 'Metropolitana'
$ws.Cells.Item(345, 5).Value = 13
$ws.Cells.Item(345, 6).Value = 100112028
$ws.Cells.Item(345, 7).Value = 'Sandia'
$ws.Cells.Item(345, 8).Value = 'Sin especificar'
$ws.Cells.Item(345, 17).Value = 1
$ws.Cells.Item(345, 18).Value = 'Hortaliza'
$ws.Cells.Item(345, 4).Value = 44525
$ws.Cells.Item(345, 9).Value = 'Segunda'
$ws.Cells.Item(345, 10).Value = 160
$ws.Cells.Item(345, 11).Value = 500
$ws.Cells.Item(345, 12).Value = 600
$ws.Cells.Item(345, 13).Value = 550
$ws.Cells.Item(345, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(345, 15).Value = 'Perú'
$ws.Cells.Item(345, 16).Value = 550
